$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 904
$ws.Range("C2").Value = 644
$ws.Range("D2").Value = 541

$ws.Range("B3").Value = 502
$ws.Range("C3").Value = 690
$ws.Range("D3").Value = 169

$ws.Range("B4").Value = 818
$ws.Range("C4").Value = 1334
$ws.Range("D4").Value = 160

$ws.Range("C5").Value = 1357
